$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 3066.639
$ws.Range("I76").Value = 3041.147
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 3041.147
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -2726.147
$ws.Range("N76").Value = -4130
# Row 79
$ws.Range("H79").Value = 3066.639
$ws.Range("I79").Value = 3041.147
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 3041.147
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -1949.147
$ws.Range("N79").Value = -5684
# Row 100
$ws.Range("H100").Value = 14287526
$ws.Range("I100").Value = 15386468
$ws.Range("J100").Value = 1276
$ws.Range("K100").Value = 15386468
$ws.Range("L100").Value = 1276
$ws.Range("M100").Value = -15385927
$ws.Range("N100").Value = -2358
# Row 127
$ws.Range("H127").Value = 1591.6578
$ws.Range("I127").Value = 820.8
$ws.Range("K127").Value = 2462.4
$ws.Range("M127").Value = 2497.6
# Row 129
$ws.Range("H129").Value = 914.77
$ws.Range("I129").Value = 475
$ws.Range("J129").Value = 933.09375
$ws.Range("K129").Value = 1425
$ws.Range("L129").Value = 2799.28125
$ws.Range("M129").Value = 3575
$ws.Range("N129").Value = -12799.28125
# Row 137
$ws.Range("H137").Value = 3533.4
$ws.Range("I137").Value = 1750.9412
$ws.Range("J137").Value = 7321.125
$ws.Range("K137").Value = 5252.8236
$ws.Range("L137").Value = 21963.375
$ws.Range("M137").Value = -2702.8236
$ws.Range("N137").Value = -27063.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5861.8545
$ws.Range("I32").Value = 4962.089
$ws.Range("K32").Value = 4962.089
$ws.Range("M32").Value = -4675.089
# Row 61
$ws.Range("H61").Value = 1642.3
$ws.Range("I61").Value = 1002.875
$ws.Range("J61").Value = 4200
$ws.Range("K61").Value = 1002.875
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -790.875
$ws.Range("N61").Value = -4624
# Row 74
$ws.Range("H74").Value = 3673.0667
$ws.Range("I74").Value = 3371.1428
$ws.Range("K74").Value = 3371.1428
$ws.Range("M74").Value = -2497.1428
# Row 77
$ws.Range("H77").Value = 3673.0667
$ws.Range("I77").Value = 3371.1428
$ws.Range("K77").Value = 16855.714
$ws.Range("M77").Value = -12487.714
# Row 136
$ws.Range("H136").Value = 1642.3
$ws.Range("I136").Value = 1002.875
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 3008.625
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -458.625
$ws.Range("N136").Value = -17700

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1741.25
$ws.Range("I86").Value = 1277.3334
$ws.Range("J86").Value = 2019.6
$ws.Range("K86").Value = 1277.3334
$ws.Range("L86").Value = 2019.6
$ws.Range("M86").Value = -154.3334
$ws.Range("N86").Value = -4265.6
# Row 89
$ws.Range("H89").Value = 1741.25
$ws.Range("I89").Value = 1277.3334
$ws.Range("J89").Value = 2019.6
$ws.Range("K89").Value = 6386.666999999999
$ws.Range("L89").Value = 10098
$ws.Range("M89").Value = -770.6669999999995
$ws.Range("N89").Value = -21330
# Row 94
$ws.Range("H94").Value = 1073.4286
$ws.Range("I94").Value = 1078.3077
$ws.Range("K94").Value = 1078.3077
$ws.Range("M94").Value = -627.3077000000001
# Row 99
$ws.Range("H99").Value = 1843.0714
$ws.Range("I99").Value = 1075.5
$ws.Range("J99").Value = 2866.5
$ws.Range("K99").Value = 1075.5
$ws.Range("L99").Value = 2866.5
$ws.Range("M99").Value = 422.5
$ws.Range("N99").Value = -5862.5
# Row 134
$ws.Range("H134").Value = 1943.8857
$ws.Range("I134").Value = 1319.3478
$ws.Range("J134").Value = 3140.9167
$ws.Range("K134").Value = 3958.0434
$ws.Range("L134").Value = 9422.750100000001
$ws.Range("M134").Value = -1423.0434
$ws.Range("N134").Value = -14492.7501

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 12197361
$ws.Range("I31").Value = 981.6667
$ws.Range("K31").Value = 981.6667
$ws.Range("M31").Value = -686.6667
# Row 34
$ws.Range("H34").Value = 12197361
$ws.Range("I34").Value = 981.6667
$ws.Range("K34").Value = 981.6667
$ws.Range("M34").Value = -779.6667
# Row 58
$ws.Range("H58").Value = 1845.0541
$ws.Range("I58").Value = 1631.7258
$ws.Range("J58").Value = 2947.25
$ws.Range("K58").Value = 1631.7258
$ws.Range("L58").Value = 2947.25
$ws.Range("M58").Value = -1428.7258
$ws.Range("N58").Value = -3353.25
# Row 105
$ws.Range("H105").Value = 1518.1111
$ws.Range("I105").Value = 1511.4839
$ws.Range("J105").Value = 1559.2
$ws.Range("K105").Value = 1511.4839
$ws.Range("L105").Value = 1559.2
$ws.Range("M105").Value = 235.5161000000001
$ws.Range("N105").Value = -5053.2
# Row 132
$ws.Range("H132").Value = 2419.1794
$ws.Range("I132").Value = 1537.2759
$ws.Range("J132").Value = 4976.7
$ws.Range("K132").Value = 4611.8277
$ws.Range("L132").Value = 14930.1
$ws.Range("M132").Value = -2081.8277
$ws.Range("N132").Value = -19990.1
# Row 134
$ws.Range("H134").Value = 6642.3335
$ws.Range("I134").Value = 8499.691999999999
$ws.Range("J134").Value = 3624.125
$ws.Range("K134").Value = 25499.076
$ws.Range("L134").Value = 10872.375
$ws.Range("M134").Value = -22964.076
$ws.Range("N134").Value = -15942.375
# Row 136
$ws.Range("H136").Value = 1845.0541
$ws.Range("I136").Value = 1631.7258
$ws.Range("J136").Value = 2947.25
$ws.Range("K136").Value = 4895.1774
$ws.Range("L136").Value = 8841.75
$ws.Range("M136").Value = -2345.1774
$ws.Range("N136").Value = -13941.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 727.8788
$ws.Range("I113").Value = 595.4737
$ws.Range("J113").Value = 907.5714
$ws.Range("K113").Value = 1786.4211
$ws.Range("L113").Value = 2722.7142
$ws.Range("M113").Value = 383.5789
$ws.Range("N113").Value = -7062.7142
# Row 131
$ws.Range("H131").Value = 9434857
$ws.Range("I131").Value = 100000270
$ws.Range("J131").Value = 959.6042
$ws.Range("K131").Value = 300000810
$ws.Range("L131").Value = 2878.8126
$ws.Range("M131").Value = -299995770
$ws.Range("N131").Value = -12958.8126

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 104
$ws.Range("H104").Value = 31500
$ws.Range("J104").Value = 31500
$ws.Range("L104").Value = 31500
$ws.Range("N104").Value = -38488
# Row 123
$ws.Range("H123").Value = 11004.467
$ws.Range("J123").Value = 11004.467
$ws.Range("L123").Value = 11004.467
$ws.Range("N123").Value = -15904.467

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 8457.237999999999
$ws.Range("I40").Value = 9979.857
$ws.Range("J40").Value = 7695.9287
$ws.Range("K40").Value = 9979.857
$ws.Range("L40").Value = 7695.9287
$ws.Range("M40").Value = -9843.857
$ws.Range("N40").Value = -7967.9287
# Row 93
$ws.Range("H93").Value = 7938577
$ws.Range("I93").Value = 12347289
$ws.Range("J93").Value = 2895.8
$ws.Range("K93").Value = 12347289
$ws.Range("L93").Value = 2895.8
$ws.Range("M93").Value = -12346041
$ws.Range("N93").Value = -5391.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 29450
$ws.Range("J16").Value = 29450
$ws.Range("L16").Value = 29450
$ws.Range("N16").Value = -30034
# Row 62
$ws.Range("H62").Value = 35333
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 51500
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 51500
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -52748
# Row 65
$ws.Range("H65").Value = 35333
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 51500
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 257500
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -263740
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 122
$ws.Range("H122").Value = 2734
$ws.Range("I122").Value = 1644.95
$ws.Range("J122").Value = 6364.1665
$ws.Range("K122").Value = 4934.85
$ws.Range("L122").Value = 19092.4995
$ws.Range("M122").Value = -2484.85
$ws.Range("N122").Value = -23992.4995
# Row 132
$ws.Range("H132").Value = 6290674.5
$ws.Range("I132").Value = 853.6
$ws.Range("K132").Value = 2560.8
$ws.Range("M132").Value = -30.80000000000018
# Row 136
$ws.Range("H136").Value = 2392.075
$ws.Range("I136").Value = 747.5185
$ws.Range("J136").Value = 5807.6924
$ws.Range("K136").Value = 2242.5555
$ws.Range("L136").Value = 17423.0772
$ws.Range("M136").Value = 307.4445000000001
$ws.Range("N136").Value = -22523.0772

